$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header renames: hire_date -> start_date, active -> is_active ---
$ws.Range("I1").Value = "start_date"
$ws.Range("K1").Value = "is_active"

# --- New employee row 14 (id 13): Axel ---
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "Axel"
$ws.Range("C14").Value = "Axel"
$ws.Range("F14").Value = "Front"
$ws.Range("I2").Copy()
$ws.Range("I14").PasteSpecial(-4122)
$ws.Range("I14").Value = 44866
$ws.Range("I2").Copy()
$ws.Range("J14").PasteSpecial(-4122)
$ws.Range("J14").Value = 45291
$ws.Range("K14").Value = "No"
$ws.Range("L14").Value = "Line Cook"

# --- New employee row 15 (id 14): Thammanit "Hui" Sawangchad ---
$ws.Range("A15").Value = 14
$ws.Range("C15").Value = "Hui"
$ws.Range("B15").Value = "Thammanit"
$ws.Range("E15").Value = "Sawangchad"
$ws.Range("F15").Value = "Back"
$ws.Range("I2").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$ws.Range("I15").Value = 45397
$ws.Range("I2").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 45504
$ws.Range("K15").Value = "No"
$ws.Range("L15").Value = "Line Cook"
$ws.Range("M15").Value = "10993 Scarlet St, Loma Linda, 92354"

# --- Updated address for rows 2 & 3 (company moved to a new address) ---
$ws.Range("M2").Value = "25581 Nikcs Avenue, Loma Linda, 92354"
$ws.Range("M3").Value = "25581 Nikcs Avenue, Loma Linda, 92354"

# --- Leave selection where the author left it when saving ---
$null = $ws.Range("M24").Select()
